{"js": "// Update the date line and the twenty-five \"two-digit \u00f7 one-digit\" problems\n// to the next day's worksheet values. Replacements are applied in document\n// order (paragraph order), matching the ordered list below exactly so that\n// repeated values (e.g. \"94\u00f77=\" used twice as a target) are handled safely.\nconst replacements = [\n  [\"2024-11-20 Wednesday\", \"2024-11-21 Thursday\"],\n  [\"76\u00f76=\", \"79\u00f73=\"],\n  [\"15\u00f78=\", \"69\u00f78=\"],\n  [\"20\u00f78=\", \"30\u00f79=\"],\n  [\"82\u00f74=\", \"46\u00f73=\"],\n  [\"52\u00f74=\", \"71\u00f72=\"],\n  [\"61\u00f72=\", \"92\u00f78=\"],\n  [\"29\u00f79=\", \"83\u00f79=\"],\n  [\"94\u00f77=\", \"53\u00f77=\"],\n  [\"45\u00f76=\", \"35\u00f72=\"],\n  [\"43\u00f73=\", \"75\u00f76=\"],\n  [\"88\u00f75=\", \"94\u00f77=\"],\n  [\"47\u00f72=\", \"45\u00f78=\"],\n  [\"87\u00f73=\", \"45\u00f77=\"],\n  [\"13\u00f79=\", \"70\u00f72=\"],\n  [\"41\u00f79=\", \"72\u00f78=\"],\n  [\"25\u00f73=\", \"46\u00f76=\"],\n  [\"39\u00f78=\", \"38\u00f79=\"],\n  [\"46\u00f78=\", \"95\u00f76=\"],\n  [\"27\u00f78=\", \"51\u00f74=\"],\n  [\"56\u00f75=\", \"95\u00f79=\"],\n  [\"51\u00f75=\", \"92\u00f79=\"],\n  [\"72\u00f77=\", \"36\u00f74=\"],\n  [\"76\u00f73=\", \"50\u00f74=\"],\n  [\"17\u00f79=\", \"64\u00f77=\"],\n  [\"55\u00f77=\", \"94\u00f77=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  if (text === \"\") {\n    continue;\n  }\n  if (idx >= replacements.length) {\n    break;\n  }\n  const [expected, next] = replacements[idx];\n  if (text !== expected) {\n    throw new Error(\n      `Unexpected paragraph text at position ${idx}: expected \"${expected}\" but found \"${text}\"`\n    );\n  }\n  p.insertText(next, Word.InsertLocation.replace);\n  idx++;\n}\n\nawait context.sync();\n\nif (idx !== replacements.length) {\n  throw new Error(`Only applied ${idx} of ${replacements.length} replacements`);\n}\n", "ps1": "# Update the date line and the twenty-five \"two-digit / one-digit\" division\n# problems to the next day's worksheet values. Replacements are applied in\n# document order (paragraph order) so repeated values (e.g. \"94\u00f77=\" used as\n# a target twice) are handled unambiguously.\n$replacements = @(\n    @(\"2024-11-20 Wednesday\", \"2024-11-21 Thursday\"),\n    @(\"76\u00f76=\", \"79\u00f73=\"),\n    @(\"15\u00f78=\", \"69\u00f78=\"),\n    @(\"20\u00f78=\", \"30\u00f79=\"),\n    @(\"82\u00f74=\", \"46\u00f73=\"),\n    @(\"52\u00f74=\", \"71\u00f72=\"),\n    @(\"61\u00f72=\", \"92\u00f78=\"),\n    @(\"29\u00f79=\", \"83\u00f79=\"),\n    @(\"94\u00f77=\", \"53\u00f77=\"),\n    @(\"45\u00f76=\", \"35\u00f72=\"),\n    @(\"43\u00f73=\", \"75\u00f76=\"),\n    @(\"88\u00f75=\", \"94\u00f77=\"),\n    @(\"47\u00f72=\", \"45\u00f78=\"),\n    @(\"87\u00f73=\", \"45\u00f77=\"),\n    @(\"13\u00f79=\", \"70\u00f72=\"),\n    @(\"41\u00f79=\", \"72\u00f78=\"),\n    @(\"25\u00f73=\", \"46\u00f76=\"),\n    @(\"39\u00f78=\", \"38\u00f79=\"),\n    @(\"46\u00f78=\", \"95\u00f76=\"),\n    @(\"27\u00f78=\", \"51\u00f74=\"),\n    @(\"56\u00f75=\", \"95\u00f79=\"),\n    @(\"51\u00f75=\", \"92\u00f79=\"),\n    @(\"72\u00f77=\", \"36\u00f74=\"),\n    @(\"76\u00f73=\", \"50\u00f74=\"),\n    @(\"17\u00f79=\", \"64\u00f77=\"),\n    @(\"55\u00f77=\", \"94\u00f77=\")\n)\n\n$d = $word.ActiveDocument\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $range = $p.Range\n    # Paragraph/cell-ending marks (CR and, inside table cells, the cell mark)\n    # are included in Range.Text; strip them before comparing the content.\n    $clean = $range.Text.TrimEnd([char]13, [char]7)\n\n    if ($clean.Length -eq 0) {\n        continue\n    }\n\n    if ($idx -ge $replacements.Count) {\n        break\n    }\n\n    $expected = $replacements[$idx][0]\n    $next = $replacements[$idx][1]\n\n    if ($clean -ne $expected) {\n        throw \"Unexpected paragraph text at position $idx`: expected '$expected' but found '$clean'\"\n    }\n\n    # Assigning Range.Text replaces only the content, preserving the\n    # paragraph/cell mark and the run's existing character formatting.\n    $range.Text = $next\n    $idx++\n}\n\nif ($idx -ne $replacements.Count) {\n    throw \"Only applied $idx of $($replacements.Count) replacements\"\n}\n"}
